# Insert a new price record as row 59 on the active sheet, pushing the
# existing rows 59..76 down to 60..77 (and extending the used range /
# dimension to A1:R77 accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59:76 down by one, creating a blank row 59.
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new "Española" record.
$ws.Cells.Item(59, 1).Value2  = 7
$ws.Cells.Item(59, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value2  = "Ñuble"
$ws.Cells.Item(59, 4).Value2  = 44855
$ws.Cells.Item(59, 5).Value2  = 16
$ws.Cells.Item(59, 6).Value2  = 100112013
$ws.Cells.Item(59, 7).Value2  = "Alcachofa"
$ws.Cells.Item(59, 8).Value2  = "Española"
$ws.Cells.Item(59, 9).Value2  = "Primera"
$ws.Cells.Item(59, 10).Value2 = 60
$ws.Cells.Item(59, 11).Value2 = 12000
$ws.Cells.Item(59, 12).Value2 = 12000
$ws.Cells.Item(59, 13).Value2 = 12000
$ws.Cells.Item(59, 14).Value2 = "$/caja 30 unidades"
$ws.Cells.Item(59, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(59, 16).Value2 = 400
$ws.Cells.Item(59, 17).Value2 = 30
$ws.Cells.Item(59, 18).Value2 = "Hortaliza"
